$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2264.5217
$ws.Range("I40").Value = 2167.2666
$ws.Range("J40").Value = 2446.875
$ws.Range("K40").Value = 2167.2666
$ws.Range("L40").Value = 2446.875
$ws.Range("M40").Value = -1992.2666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 310.9375
$ws.Range("I53").Value = 234
$ws.Range("J53").Value = 439.16666
$ws.Range("K53").Value = 234
$ws.Range("L53").Value = 439.16666
$ws.Range("M53").Value = 403

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5463.467
$ws.Range("I61").Value = 1660.5
$ws.Range("J61").Value = 7998.778
$ws.Range("K61").Value = 1660.5
$ws.Range("L61").Value = 7998.778
$ws.Range("M61").Value = -1448.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1917.75
$ws.Range("I110").Value = 1917.75
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1917.75
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 127.25
$ws.Range("N110").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2279.1
$ws.Range("I132").Value = 2224.125
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 6672.375
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -4142.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5463.467
$ws.Range("I136").Value = 1660.5
$ws.Range("J136").Value = 7998.778
$ws.Range("K136").Value = 4981.5
$ws.Range("L136").Value = 23996.334
$ws.Range("M136").Value = -2431.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1032.8
$ws.Range("I20").Value = 1162.6666
$ws.Range("J20").Value = 838
$ws.Range("K20").Value = 1162.6666
$ws.Range("L20").Value = 838
$ws.Range("M20").Value = -915.6666
$ws.Range("N20").Value = -1332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 810.6429000000001
$ws.Range("I22").Value = 606.9
$ws.Range("J22").Value = 1320
$ws.Range("K22").Value = 606.9
$ws.Range("L22").Value = 1320
$ws.Range("M22").Value = -433.9
$ws.Range("N22").Value = -1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7159.7
$ws.Range("I94").Value = 7833.1113
$ws.Range("J94").Value = 1099
$ws.Range("K94").Value = 7833.1113
$ws.Range("L94").Value = 1099
$ws.Range("M94").Value = -7382.1113

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1567.95
$ws.Range("I16").Value = 1083.2667
$ws.Range("J16").Value = 3022
$ws.Range("K16").Value = 1083.2667
$ws.Range("L16").Value = 3022
$ws.Range("M16").Value = -796.2666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2176
$ws.Range("I58").Value = 2548.8
$ws.Range("J58").Value = 1909.7142
$ws.Range("K58").Value = 2548.8
$ws.Range("L58").Value = 1909.7142
$ws.Range("M58").Value = -2345.8
$ws.Range("N58").Value = -2315.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 851.75
$ws.Range("I107").Value = 472.85715
$ws.Range("J107").Value = 1382.2
$ws.Range("K107").Value = 472.85715
$ws.Range("L107").Value = 1382.2
$ws.Range("M107").Value = 1447.14285
$ws.Range("N107").Value = -5222.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1567.95
$ws.Range("I113").Value = 1083.2667
$ws.Range("J113").Value = 3022
$ws.Range("K113").Value = 1083.2667
$ws.Range("L113").Value = 3022
$ws.Range("M113").Value = 1086.7333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2176
$ws.Range("I136").Value = 2548.8
$ws.Range("J136").Value = 1909.7142
$ws.Range("K136").Value = 7646.400000000001
$ws.Range("L136").Value = 5729.142599999999
$ws.Range("M136").Value = -5096.400000000001
$ws.Range("N136").Value = -10829.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 19920728
$ws.Range("I4").Value = 22000112
$ws.Range("J4").Value = 9523809
$ws.Range("K4").Value = 66000336
$ws.Range("L4").Value = 28571427
$ws.Range("M4").Value = -66000224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2192.75
$ws.Range("I5").Value = 3248.6
$ws.Range("J5").Value = 433
$ws.Range("K5").Value = 9745.799999999999
$ws.Range("L5").Value = 1299
$ws.Range("M5").Value = -9633.799999999999
$ws.Range("N5").Value = -1523

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 514.7778
$ws.Range("I92").Value = 385.27777
$ws.Range("J92").Value = 773.7778
$ws.Range("K92").Value = 1155.83331
$ws.Range("L92").Value = 2321.3334
$ws.Range("M92").Value = 92.16669000000002
$ws.Range("N92").Value = -4817.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 217.16667
$ws.Range("I98").Value = 174
$ws.Range("J98").Value = 303.5
$ws.Range("K98").Value = 522
$ws.Range("L98").Value = 910.5
$ws.Range("M98").Value = 976

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1221.5454
$ws.Range("I113").Value = 1559.5
$ws.Range("J113").Value = 1146.4445
$ws.Range("K113").Value = 4678.5
$ws.Range("L113").Value = 3439.3335
$ws.Range("M113").Value = -2508.5
$ws.Range("N113").Value = -7779.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2646.6875
$ws.Range("I131").Value = 1437.25
$ws.Range("J131").Value = 3049.8333
$ws.Range("K131").Value = 4311.75
$ws.Range("L131").Value = 9149.499899999999
$ws.Range("M131").Value = 728.25
$ws.Range("N131").Value = -19229.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2192.75
$ws.Range("I135").Value = 3248.6
$ws.Range("J135").Value = 433
$ws.Range("K135").Value = 29237.4
$ws.Range("L135").Value = 3897
$ws.Range("M135").Value = -26702.4
$ws.Range("N135").Value = -8967

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2126.3333
$ws.Range("I80").Value = 2126.3333
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2126.3333
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1128.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2126.3333
$ws.Range("I83").Value = 2126.3333
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 10631.6665
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -5639.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 189
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 189
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 189
$ws.Range("M55").Value = ""
$ws.Range("N55").Value = -535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1929.7858
$ws.Range("I82").Value = 1461.9
$ws.Range("J82").Value = 3099.5
$ws.Range("K82").Value = 1461.9
$ws.Range("L82").Value = 3099.5
$ws.Range("M82").Value = -1100.9
$ws.Range("N82").Value = -3821.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1929.7858
$ws.Range("I85").Value = 1461.9
$ws.Range("J85").Value = 3099.5
$ws.Range("K85").Value = 1461.9
$ws.Range("L85").Value = 3099.5
$ws.Range("M85").Value = -213.9000000000001
$ws.Range("N85").Value = -5595.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4899
$ws.Range("I136").Value = 4881.1665
$ws.Range("J136").Value = 4934.6665
$ws.Range("K136").Value = 14643.4995
$ws.Range("L136").Value = 14803.9995
$ws.Range("M136").Value = -12093.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 110390
$ws.Range("I137").Value = 110390
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 110390
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -105290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 51998.332
$ws.Range("I70").Value = 51995
$ws.Range("J70").Value = 52000
$ws.Range("K70").Value = 51995
$ws.Range("L70").Value = 52000
$ws.Range("M70").Value = -51680
$ws.Range("N70").Value = -52630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 51998.332
$ws.Range("I73").Value = 51995
$ws.Range("J73").Value = 52000
$ws.Range("K73").Value = 51995
$ws.Range("L73").Value = 52000
$ws.Range("M73").Value = -50903
$ws.Range("N73").Value = -54184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 60000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 60000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 74214.25
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 74214.25
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 74214.25
$ws.Range("N140").Value = -84574.25
